$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-01 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-02 Thursday", 2) | Out-Null
$d.Content.Find.Execute("16+27=43", $true, $false, $false, $false, $false, $true, 1, $false, "64+3=67", 2) | Out-Null
$d.Content.Find.Execute("81-33=48", $true, $false, $false, $false, $false, $true, 1, $false, "23+17=40", 2) | Out-Null
$d.Content.Find.Execute("14+61=75", $true, $false, $false, $false, $false, $true, 1, $false, "99-14=85", 2) | Out-Null
$d.Content.Find.Execute("0+31=31", $true, $false, $false, $false, $false, $true, 1, $false, "96+2=98", 2) | Out-Null
$d.Content.Find.Execute("38+2=40", $true, $false, $false, $false, $false, $true, 1, $false, "71+22=93", 2) | Out-Null
$d.Content.Find.Execute("58-49=9", $true, $false, $false, $false, $false, $true, 1, $false, "40+2=42", 2) | Out-Null
$d.Content.Find.Execute("4+61=65", $true, $false, $false, $false, $false, $true, 1, $false, "3+65=68", 2) | Out-Null
$d.Content.Find.Execute("68+30=98", $true, $false, $false, $false, $false, $true, 1, $false, "74+14=88", 2) | Out-Null
$d.Content.Find.Execute("74-58=16", $true, $false, $false, $false, $false, $true, 1, $false, "68-5=63", 2) | Out-Null
$d.Content.Find.Execute("49-35=14", $true, $false, $false, $false, $false, $true, 1, $false, "97-1=96", 2) | Out-Null
$d.Content.Find.Execute("46-18=28", $true, $false, $false, $false, $false, $true, 1, $false, "18+6=24", 2) | Out-Null
$d.Content.Find.Execute("59-30=29", $true, $false, $false, $false, $false, $true, 1, $false, "39-11=28", 2) | Out-Null
$d.Content.Find.Execute("9+90=99", $true, $false, $false, $false, $false, $true, 1, $false, "5+73=78", 2) | Out-Null
$d.Content.Find.Execute("36+32=68", $true, $false, $false, $false, $false, $true, 1, $false, "78-66=12", 2) | Out-Null
$d.Content.Find.Execute("84-26=58", $true, $false, $false, $false, $false, $true, 1, $false, "20-9=11", 2) | Out-Null
$d.Content.Find.Execute("87-40=47", $true, $false, $false, $false, $false, $true, 1, $false, "21+12=33", 2) | Out-Null
$d.Content.Find.Execute("8+19=27", $true, $false, $false, $false, $false, $true, 1, $false, "81-60=21", 2) | Out-Null
$d.Content.Find.Execute("19+35=54", $true, $false, $false, $false, $false, $true, 1, $false, "95-40=55", 2) | Out-Null
$d.Content.Find.Execute("65-8=57", $true, $false, $false, $false, $false, $true, 1, $false, "58-34=24", 2) | Out-Null
$d.Content.Find.Execute("58-16=42", $true, $false, $false, $false, $false, $true, 1, $false, "92-83=9", 2) | Out-Null
$d.Content.Find.Execute("51+22=73", $true, $false, $false, $false, $false, $true, 1, $false, "92-8=84", 2) | Out-Null
$d.Content.Find.Execute("40-30=10", $true, $false, $false, $false, $false, $true, 1, $false, "54-26=28", 2) | Out-Null
$d.Content.Find.Execute("23-6=17", $true, $false, $false, $false, $false, $true, 1, $false, "97-67=30", 2) | Out-Null
$d.Content.Find.Execute("85-52=33", $true, $false, $false, $false, $false, $true, 1, $false, "7+15=22", 2) | Out-Null
$d.Content.Find.Execute("98-46=52", $true, $false, $false, $false, $false, $true, 1, $false, "94-54=40", 2) | Out-Null
$d.Content.Find.Execute("64-3=61", $true, $false, $false, $false, $false, $true, 1, $false, "77+2=79", 2) | Out-Null
$d.Content.Find.Execute("38+55=93", $true, $false, $false, $false, $false, $true, 1, $false, "75-51=24", 2) | Out-Null
$d.Content.Find.Execute("84+11=95", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=20", 2) | Out-Null
$d.Content.Find.Execute("42+10=52", $true, $false, $false, $false, $false, $true, 1, $false, "53-4=49", 2) | Out-Null
$d.Content.Find.Execute("82-76=6", $true, $false, $false, $false, $false, $true, 1, $false, "66-10=56", 2) | Out-Null
$d.Content.Find.Execute("6+25=31", $true, $false, $false, $false, $false, $true, 1, $false, "85-43=42", 2) | Out-Null
$d.Content.Find.Execute("67-51=16", $true, $false, $false, $false, $false, $true, 1, $false, "37+19=56", 2) | Out-Null
$d.Content.Find.Execute("29+66=95", $true, $false, $false, $false, $false, $true, 1, $false, "3+79=82", 2) | Out-Null
$d.Content.Find.Execute("18+1=19", $true, $false, $false, $false, $false, $true, 1, $false, "78-9=69", 2) | Out-Null
$d.Content.Find.Execute("88-64=24", $true, $false, $false, $false, $false, $true, 1, $false, "98-94=4", 2) | Out-Null
$d.Content.Find.Execute("8+88=96", $true, $false, $false, $false, $false, $true, 1, $false, "50+48=98", 2) | Out-Null
$d.Content.Find.Execute("38+17=55", $true, $false, $false, $false, $false, $true, 1, $false, "76+9=85", 2) | Out-Null
$d.Content.Find.Execute("37+38=75", $true, $false, $false, $false, $false, $true, 1, $false, "52-35=17", 2) | Out-Null
$d.Content.Find.Execute("13+7=20", $true, $false, $false, $false, $false, $true, 1, $false, "75+22=97", 2) | Out-Null
$d.Content.Find.Execute("99-22=77", $true, $false, $false, $false, $false, $true, 1, $false, "20+16=36", 2) | Out-Null
$d.Content.Find.Execute("29+24=53", $true, $false, $false, $false, $false, $true, 1, $false, "38-28=10", 2) | Out-Null
$d.Content.Find.Execute("26+10=36", $true, $false, $false, $false, $false, $true, 1, $false, "76+19=95", 2) | Out-Null
$d.Content.Find.Execute("0+80=80", $true, $false, $false, $false, $false, $true, 1, $false, "83-62=21", 2) | Out-Null
$d.Content.Find.Execute("18+10=28", $true, $false, $false, $false, $false, $true, 1, $false, "86-55=31", 2) | Out-Null
$d.Content.Find.Execute("0+12=12", $true, $false, $false, $false, $false, $true, 1, $false, "21+58=79", 2) | Out-Null
$d.Content.Find.Execute("66-44=22", $true, $false, $false, $false, $false, $true, 1, $false, "1+8=9", 2) | Out-Null
$d.Content.Find.Execute("34-26=8", $true, $false, $false, $false, $false, $true, 1, $false, "87+5=92", 2) | Out-Null
$d.Content.Find.Execute("25+71=96", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=65", 2) | Out-Null
$d.Content.Find.Execute("37+60=97", $true, $false, $false, $false, $false, $true, 1, $false, "56-41=15", 2) | Out-Null
$d.Content.Find.Execute("62+5=67", $true, $false, $false, $false, $false, $true, 1, $false, "68-40=28", 2) | Out-Null
$d.Content.Find.Execute("41+10=51", $true, $false, $false, $false, $false, $true, 1, $false, "14+67=81", 2) | Out-Null
$d.Content.Find.Execute("69+24=93", $true, $false, $false, $false, $false, $true, 1, $false, "28+60=88", 2) | Out-Null
$d.Content.Find.Execute("5+46=51", $true, $false, $false, $false, $false, $true, 1, $false, "38+3=41", 2) | Out-Null
$d.Content.Find.Execute("70-69=1", $true, $false, $false, $false, $false, $true, 1, $false, "26+39=65", 2) | Out-Null
$d.Content.Find.Execute("53-17=36", $true, $false, $false, $false, $false, $true, 1, $false, "34+58=92", 2) | Out-Null
$d.Content.Find.Execute("51+7=58", $true, $false, $false, $false, $false, $true, 1, $false, "89-65=24", 2) | Out-Null
$d.Content.Find.Execute("77+18=95", $true, $false, $false, $false, $false, $true, 1, $false, "68-14=54", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $true, $false, $false, $false, $false, $true, 1, $false, "0+81=81", 2) | Out-Null
$d.Content.Find.Execute("79+20=99", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=37", 2) | Out-Null
$d.Content.Find.Execute("78-51=27", $true, $false, $false, $false, $false, $true, 1, $false, "15+61=76", 2) | Out-Null
$d.Content.Find.Execute("39-0=39", $true, $false, $false, $false, $false, $true, 1, $false, "24+17=41", 2) | Out-Null
$d.Content.Find.Execute("64+8=72", $true, $false, $false, $false, $false, $true, 1, $false, "93-87=6", 2) | Out-Null
$d.Content.Find.Execute("45-4=41", $true, $false, $false, $false, $false, $true, 1, $false, "66+26=92", 2) | Out-Null
$d.Content.Find.Execute("49-39=10", $true, $false, $false, $false, $false, $true, 1, $false, "14-4=10", 2) | Out-Null
$d.Content.Find.Execute("43-26=17", $true, $false, $false, $false, $false, $true, 1, $false, "36+46=82", 2) | Out-Null
$d.Content.Find.Execute("74-70=4", $true, $false, $false, $false, $false, $true, 1, $false, "1+8=9", 2) | Out-Null
$d.Content.Find.Execute("94-58=36", $true, $false, $false, $false, $false, $true, 1, $false, "42-22=20", 2) | Out-Null
$d.Content.Find.Execute("88+4=92", $true, $false, $false, $false, $false, $true, 1, $false, "22+46=68", 2) | Out-Null
$d.Content.Find.Execute("8+60=68", $true, $false, $false, $false, $false, $true, 1, $false, "20+9=29", 2) | Out-Null
$d.Content.Find.Execute("88+9=97", $true, $false, $false, $false, $false, $true, 1, $false, "8+72=80", 2) | Out-Null
$d.Content.Find.Execute("42+37=79", $true, $false, $false, $false, $false, $true, 1, $false, "22+46=68", 2) | Out-Null
$d.Content.Find.Execute("25+67=92", $true, $false, $false, $false, $false, $true, 1, $false, "83-41=42", 2) | Out-Null
$d.Content.Find.Execute("47+6=53", $true, $false, $false, $false, $false, $true, 1, $false, "16+52=68", 2) | Out-Null
$d.Content.Find.Execute("31-16=15", $true, $false, $false, $false, $false, $true, 1, $false, "74-9=65", 2) | Out-Null
$d.Content.Find.Execute("79-37=42", $true, $false, $false, $false, $false, $true, 1, $false, "62+9=71", 2) | Out-Null
$d.Content.Find.Execute("98-8=90", $true, $false, $false, $false, $false, $true, 1, $false, "54+20=74", 2) | Out-Null
$d.Content.Find.Execute("71-28=43", $true, $false, $false, $false, $false, $true, 1, $false, "32+66=98", 2) | Out-Null
$d.Content.Find.Execute("91-30=61", $true, $false, $false, $false, $false, $true, 1, $false, "31+35=66", 2) | Out-Null
$d.Content.Find.Execute("73-42=31", $true, $false, $false, $false, $false, $true, 1, $false, "63+33=96", 2) | Out-Null
$d.Content.Find.Execute("72-66=6", $true, $false, $false, $false, $false, $true, 1, $false, "35-33=2", 2) | Out-Null
$d.Content.Find.Execute("91-49=42", $true, $false, $false, $false, $false, $true, 1, $false, "33-18=15", 2) | Out-Null
$d.Content.Find.Execute("83-16=67", $true, $false, $false, $false, $false, $true, 1, $false, "18+39=57", 2) | Out-Null
$d.Content.Find.Execute("86-5=81", $true, $false, $false, $false, $false, $true, 1, $false, "15-13=2", 2) | Out-Null
$d.Content.Find.Execute("68-32=36", $true, $false, $false, $false, $false, $true, 1, $false, "57+29=86", 2) | Out-Null
$d.Content.Find.Execute("31-17=14", $true, $false, $false, $false, $false, $true, 1, $false, "57+14=71", 2) | Out-Null
$d.Content.Find.Execute("30+8=38", $true, $false, $false, $false, $false, $true, 1, $false, "17+74=91", 2) | Out-Null
$d.Content.Find.Execute("24+42=66", $true, $false, $false, $false, $false, $true, 1, $false, "99-0=99", 2) | Out-Null
$d.Content.Find.Execute("90-7=83", $true, $false, $false, $false, $false, $true, 1, $false, "58-24=34", 2) | Out-Null
$d.Content.Find.Execute("8+12=20", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=93", 2) | Out-Null
$d.Content.Find.Execute("36-32=4", $true, $false, $false, $false, $false, $true, 1, $false, "12+82=94", 2) | Out-Null
$d.Content.Find.Execute("14+41=55", $true, $false, $false, $false, $false, $true, 1, $false, "98-5=93", 2) | Out-Null
$d.Content.Find.Execute("92-66=26", $true, $false, $false, $false, $false, $true, 1, $false, "84+15=99", 2) | Out-Null
$d.Content.Find.Execute("80-23=57", $true, $false, $false, $false, $false, $true, 1, $false, "50+7=57", 2) | Out-Null
$d.Content.Find.Execute("24-19=5", $true, $false, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("21-21=0", $true, $false, $false, $false, $false, $true, 1, $false, "33-6=27", 2) | Out-Null
$d.Content.Find.Execute("64-14=50", $true, $false, $false, $false, $false, $true, 1, $false, "24+6=30", 2) | Out-Null
$d.Content.Find.Execute("28-0=28", $true, $false, $false, $false, $false, $true, 1, $false, "1+69=70", 2) | Out-Null
$d.Content.Find.Execute("96-29=67", $true, $false, $false, $false, $false, $true, 1, $false, "74+17=91", 2) | Out-Null
$d.Content.Find.Execute("59+19=78", $true, $false, $false, $false, $false, $true, 1, $false, "56+17=73", 2) | Out-Null
$d.Content.Find.Execute("45-20=25", $true, $false, $false, $false, $false, $true, 1, $false, "38-17=21", 2) | Out-Null
